$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 16:05"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1688739
$ws.Range("C4").Value = 2303
$ws.Range("E4").Value = 1137638
$ws.Range("G4").Value = 56
$ws.Range("H4").Value = 99356

# Alemania (row 11)
$ws.Range("B11").Value = 180505
$ws.Range("C11").Value = 177
$ws.Range("E11").Value = 10923
$ws.Range("G11").Value = 11
$ws.Range("H11").Value = 8382

# Suiza (row 31)
$ws.Range("E31").Value = 737
$ws.Range("G31").Value = 3
$ws.Range("H31").Value = 1909

# Kenia (row 102)
$ws.Range("B102").Value = 1286
$ws.Range("C102").Value = 72
$ws.Range("D102").Value = 392
$ws.Range("E102").Value = 842
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 52

# Sri Lanka (row 103)
$ws.Range("B103").Value = 1182
$ws.Range("C103").Value = 41
$ws.Range("E103").Value = 477

# Principado de Andorra (row 121)
$ws.Range("B121").Value = 763
$ws.Range("C121").Value = 1
$ws.Range("D121").Value = 663
$ws.Range("E121").Value = 49
